$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ETLE")
$ws.Range("B2").Value = -2
